$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 188 (shifts existing rows 188-210 down to 189-211)
$ws.Rows(188).Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Range("A188").Value = 7
$ws.Range("B188").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C188").Value = "Ñuble"
$ws.Range("D188").Value = 45180
$ws.Range("E188").Value = 16
$ws.Range("F188").Value = 100112037
$ws.Range("G188").Value = "Cebollín"
$ws.Range("H188").Value = "Sin especificar"
$ws.Range("I188").Value = "Primera"
$ws.Range("J188").Value = 120
$ws.Range("K188").Value = 6000
$ws.Range("L188").Value = 6000
$ws.Range("M188").Value = 6000
$ws.Range("N188").Value = "$/paquete 36 unidades"
$ws.Range("O188").Value = "Provincia de Diguillín"
$ws.Range("P188").Value = 167
$ws.Range("Q188").Value = 36
$ws.Range("R188").Value = "Hortaliza"
